$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 29 / 31: re-point A29 to "Refactoring - duration" and A31 to "Everything wrong"
$ws.Range("A29").Value = "Refactoring - duration"
$ws.Range("A31").Value = "Everything wrong"

# New rows 33-40
$ws.Range("A33").Value = "Refactoring - num errors introduced"

$ws.Range("A34").Value = "Everything right"
$ws.Range("B34").Value = 3

$ws.Range("A35").Value = "6 errors, no coding standard, no refactoring tool (everything wrong)"
$ws.Range("B35").Value = 50

$ws.Range("A36").Value = "0 errors, no coding standard, no refactoring tool"
$ws.Range("B36").Value = 47

$ws.Range("A37").Value = "6 errors, no refactoring tool"
$ws.Range("B37").Value = 21

$ws.Range("A38").Value = "6 errors, no coding standard"
$ws.Range("B38").Value = 21

$ws.Range("A39").Value = "0 errors, no refactoring tool"
$ws.Range("B39").Value = 20

$ws.Range("A40").Value = "0 errors, no coding standard"
$ws.Range("B40").Value = 20

# Column A width update (62 "visual" chars; the host's px-quantization needs a
# nudge above 61 to land exactly on width=62 once re-serialized)
$ws.Columns.Item(1).ColumnWidth = 61.15

# Sheet view changes: scroll so row 18 is at the top, final selection B41
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("B41").Select()
